$d = $word.ActiveDocument

# 1. "update_counter()" re-ran (same value, 554) - re-apply the text via
#    Find/Replace so Word recombines the "REMAINING HOURS: " / "554" runs
#    into a single run, matching what happens when the counter routine
#    rewrites that phrase.
$d.Content.Find.Execute(
    "REMAINING HOURS: 554", $true, $false, $false, $false, $false,
    $true, 1, $false, "REMAINING HOURS: 554", 2) | Out-Null

# 2. Add the extra journal sentences right after "...with Laravel." and
#    before the two trailing line breaks.
$rng = $d.Content
$rng.Find.Execute(
    "with Laravel.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" Here, I learned the general, theoretical idea of how these files should interact with each other, along with a basic file structure within the editor. I am familiar with OOP from my time in Python and Java, but the PHP syntax is tripping me up.")
